# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet right before "总计" (mirrors the other
#    quarterly fund-holding sheets: 基金代码/基金名称/基金规模/股票总仓位/
#    仓位占比/持有市值(亿元)/仓位排名).
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet (日期/持有数量(只)/
#    持有市值(亿元)).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# NOTE: worksheet handles obtained via Worksheets.Item(...) resolve by
# *position*, so a handle fetched before an insert/add can silently start
# pointing at a different sheet once sheets shift around. Always re-fetch
# "总计" by name right before using it, and grab "2022-Q1" by name too
# (instead of trusting the reference returned by Worksheets.Add()).
# ---------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")
$styleDonor = $wb.Worksheets.Item("2021-Q4")

$null = $wb.Worksheets.Add($totalSheet)
$wb.ActiveSheet.Name = "2022-Q1"
$newSheet = $wb.Worksheets.Item("2022-Q1")

# ---- Header row (copy the bold/bordered/centered style used by the other
#      quarter sheets onto our new header cells) ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$styleDonor.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# ---- Row-index column (A) — reuse the donor's styled index-column cell ----
$styleDonor.Range("A2").Copy()
$newSheet.Range("A2:A13").PasteSpecial(-4122)

# Columns D:G hold numeric-looking figures that are stored as literal text
# in this workbook (e.g. "141.91"), same as every other quarter sheet.
# Force text entry, then drop the number-format style again afterwards so
# the cells end up format-less like the source, while staying text typed.
$newSheet.Range("B2:G13").NumberFormat = "@"

$data = @(
    @(0,  "008969", "睿远均衡价值三年持有期混合A",       "141.91", "93.01", "3.73", "5.2932", 6),
    @(1,  "000762", "汇添富绝对收益策略定期开放混合A",   "265.54", "56.00", "0.83", "2.2040", 10),
    @(2,  "519018", "汇添富均衡增长混合",                 "39.45",  "93.81", "5.40", "2.1303", 3),
    @(3,  "011271", "汇添富价值成长均衡投资混合A",       "24.78",  "92.58", "3.74", "0.9268", 9),
    @(4,  "008970", "睿远均衡价值三年持有期混合C",       "15.26",  "93.01", "3.73", "0.5692", 6),
    @(5,  "470008", "汇添富策略回报混合",                 "12.75",  "94.01", "4.07", "0.5189", 8),
    @(6,  "001166", "建信环保产业股票",                   "9.21",   "86.36", "3.80", "0.3500", 5),
    @(7,  "008140", "汇添富绝对收益策略定期开放混合C",   "22.23",  "56.00", "0.83", "0.1845", 10),
    @(8,  "010447", "中邮未来成长混合A",                 "0.44",   "79.42", "3.34", "0.0147", 10),
    @(9,  "011272", "汇添富价值成长均衡投资混合C",       "0.32",   "92.58", "3.74", "0.0120", 9),
    @(10, "001899", "东海中证社会发展安全产业主题指数",  "0.21",   "90.30", "2.24", "0.0047", 6),
    @(11, "010448", "中邮未来成长混合C",                 "0.07",   "79.42", "3.34", "0.0023", 10)
)

$row = 2
foreach ($rec in $data) {
    $newSheet.Cells.Item($row, 1).Value = $rec[0]
    $newSheet.Cells.Item($row, 2).Value = $rec[1]
    $newSheet.Cells.Item($row, 3).Value = $rec[2]
    $newSheet.Cells.Item($row, 4).Value = $rec[3]
    $newSheet.Cells.Item($row, 5).Value = $rec[4]
    $newSheet.Cells.Item($row, 6).Value = $rec[5]
    $newSheet.Cells.Item($row, 7).Value = $rec[6]
    $newSheet.Cells.Item($row, 8).Value = $rec[7]
    $row = $row + 1
}

# Drop the "@" number-format styling again now that the text is committed,
# so these data cells end up unstyled like the rest of the workbook.
$newSheet.Range("B2:G13").Style = "Normal"

# ---------------------------------------------------------------------
# "总计" sheet: push a new "2022-Q1" row in above the existing data.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalRows = @(
    @(1, "2021-Q4", 11, 7.75),
    @(2, "2021-Q3", 2,  0.03),
    @(3, "2021-Q2", 4,  0.08),
    @(4, "2021-Q1", 6,  0.54),
    @(5, "2020-Q4", 10, 0.57)
)
foreach ($rec in $totalRows) {
    $r = $rec[0] + 2
    $totalSheet.Cells.Item($r, 1).Value = $rec[0]
    $totalSheet.Cells.Item($r, 2).Value = $rec[1]
    $totalSheet.Cells.Item($r, 3).Value = $rec[2]
    $totalSheet.Cells.Item($r, 4).Value = $rec[3]
}

# Row 7 (previously non-existent) needs the same styled index-column cell
# as the rest of column A.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)
$totalSheet.Range("A7").Value = 5

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 12
$totalSheet.Range("D2").Value = 12.21
